$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-19 Tuesday", $false, $true, $false, $false, $false, $true, 1, $false, "2024-03-20 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("45+47=92", $false, $true, $false, $false, $false, $true, 1, $false, "2+71=73", 2) | Out-Null
$d.Content.Find.Execute("1+11=12", $false, $true, $false, $false, $false, $true, 1, $false, "55-53=2", 2) | Out-Null
$d.Content.Find.Execute("54-20=34", $false, $true, $false, $false, $false, $true, 1, $false, "21+18=39", 2) | Out-Null
$d.Content.Find.Execute("98-51=47", $false, $true, $false, $false, $false, $true, 1, $false, "21-21=0", 2) | Out-Null
$d.Content.Find.Execute("87-19=68", $false, $true, $false, $false, $false, $true, 1, $false, "13+63=76", 2) | Out-Null
$d.Content.Find.Execute("69-12=57", $false, $true, $false, $false, $false, $true, 1, $false, "40+38=78", 2) | Out-Null
$d.Content.Find.Execute("91-7=84", $false, $true, $false, $false, $false, $true, 1, $false, "47+36=83", 2) | Out-Null
$d.Content.Find.Execute("95-91=4", $false, $true, $false, $false, $false, $true, 1, $false, "71-3=68", 2) | Out-Null
$d.Content.Find.Execute("9+13=22", $false, $true, $false, $false, $false, $true, 1, $false, "65-18=47", 2) | Out-Null
$d.Content.Find.Execute("96-24=72", $false, $true, $false, $false, $false, $true, 1, $false, "1+58=59", 2) | Out-Null
$d.Content.Find.Execute("43-24=19", $false, $true, $false, $false, $false, $true, 1, $false, "50-7=43", 2) | Out-Null
$d.Content.Find.Execute("1+61=62", $false, $true, $false, $false, $false, $true, 1, $false, "62-59=3", 2) | Out-Null
$d.Content.Find.Execute("84-52=32", $false, $true, $false, $false, $false, $true, 1, $false, "55-29=26", 2) | Out-Null
$d.Content.Find.Execute("66-31=35", $false, $true, $false, $false, $false, $true, 1, $false, "99-68=31", 2) | Out-Null
$d.Content.Find.Execute("46+27=73", $false, $true, $false, $false, $false, $true, 1, $false, "41+47=88", 2) | Out-Null
$d.Content.Find.Execute("84-72=12", $false, $true, $false, $false, $false, $true, 1, $false, "95-81=14", 2) | Out-Null
$d.Content.Find.Execute("36+41=77", $false, $true, $false, $false, $false, $true, 1, $false, "64+29=93", 2) | Out-Null
$d.Content.Find.Execute("48-28=20", $false, $true, $false, $false, $false, $true, 1, $false, "86-75=11", 2) | Out-Null
$d.Content.Find.Execute("4+34=38", $false, $true, $false, $false, $false, $true, 1, $false, "66-24=42", 2) | Out-Null
$d.Content.Find.Execute("10+64=74", $false, $true, $false, $false, $false, $true, 1, $false, "49-11=38", 2) | Out-Null
$d.Content.Find.Execute("34+52=86", $false, $true, $false, $false, $false, $true, 1, $false, "46-5=41", 2) | Out-Null
$d.Content.Find.Execute("14+12=26", $false, $true, $false, $false, $false, $true, 1, $false, "44+21=65", 2) | Out-Null
$d.Content.Find.Execute("2+33=35", $false, $true, $false, $false, $false, $true, 1, $false, "78+10=88", 2) | Out-Null
$d.Content.Find.Execute("70-21=49", $false, $true, $false, $false, $false, $true, 1, $false, "61-39=22", 2) | Out-Null
$d.Content.Find.Execute("34+57=91", $false, $true, $false, $false, $false, $true, 1, $false, "22+57=79", 2) | Out-Null
$d.Content.Find.Execute("97-70=27", $false, $true, $false, $false, $false, $true, 1, $false, "53-1=52", 2) | Out-Null
$d.Content.Find.Execute("43-31=12", $false, $true, $false, $false, $false, $true, 1, $false, "61-59=2", 2) | Out-Null
$d.Content.Find.Execute("69-45=24", $false, $true, $false, $false, $false, $true, 1, $false, "7+7=14", 2) | Out-Null
$d.Content.Find.Execute("85+3=88", $false, $true, $false, $false, $false, $true, 1, $false, "52-8=44", 2) | Out-Null
$d.Content.Find.Execute("54+30=84", $false, $true, $false, $false, $false, $true, 1, $false, "31+5=36", 2) | Out-Null
$d.Content.Find.Execute("48+20=68", $false, $true, $false, $false, $false, $true, 1, $false, "93-46=47", 2) | Out-Null
$d.Content.Find.Execute("74-46=28", $false, $true, $false, $false, $false, $true, 1, $false, "59-4=55", 2) | Out-Null
$d.Content.Find.Execute("1+42=43", $false, $true, $false, $false, $false, $true, 1, $false, "3+65=68", 2) | Out-Null
$d.Content.Find.Execute("88-15=73", $false, $true, $false, $false, $false, $true, 1, $false, "19+65=84", 2) | Out-Null
$d.Content.Find.Execute("66-37=29", $false, $true, $false, $false, $false, $true, 1, $false, "72+6=78", 2) | Out-Null
$d.Content.Find.Execute("14+31=45", $false, $true, $false, $false, $false, $true, 1, $false, "1+60=61", 2) | Out-Null
$d.Content.Find.Execute("42-8=34", $false, $true, $false, $false, $false, $true, 1, $false, "81-9=72", 2) | Out-Null
$d.Content.Find.Execute("42-24=18", $false, $true, $false, $false, $false, $true, 1, $false, "95-78=17", 2) | Out-Null
$d.Content.Find.Execute("39+22=61", $false, $true, $false, $false, $false, $true, 1, $false, "1+67=68", 2) | Out-Null
$d.Content.Find.Execute("5+73=78", $false, $true, $false, $false, $false, $true, 1, $false, "3+0=3", 2) | Out-Null
$d.Content.Find.Execute("55+16=71", $false, $true, $false, $false, $false, $true, 1, $false, "63-6=57", 2) | Out-Null
$d.Content.Find.Execute("22+59=81", $false, $true, $false, $false, $false, $true, 1, $false, "67-15=52", 2) | Out-Null
$d.Content.Find.Execute("71-47=24", $false, $true, $false, $false, $false, $true, 1, $false, "66+9=75", 2) | Out-Null
$d.Content.Find.Execute("28+18=46", $false, $true, $false, $false, $false, $true, 1, $false, "81-62=19", 2) | Out-Null
$d.Content.Find.Execute("89-2=87", $false, $true, $false, $false, $false, $true, 1, $false, "67+4=71", 2) | Out-Null
$d.Content.Find.Execute("3+72=75", $false, $true, $false, $false, $false, $true, 1, $false, "77-75=2", 2) | Out-Null
$d.Content.Find.Execute("63+23=86", $false, $true, $false, $false, $false, $true, 1, $false, "99-44=55", 2) | Out-Null
$d.Content.Find.Execute("70+3=73", $false, $true, $false, $false, $false, $true, 1, $false, "45+14=59", 2) | Out-Null
$d.Content.Find.Execute("84-17=67", $false, $true, $false, $false, $false, $true, 1, $false, "9+3=12", 2) | Out-Null
$d.Content.Find.Execute("19+7=26", $false, $true, $false, $false, $false, $true, 1, $false, "61-28=33", 2) | Out-Null
$d.Content.Find.Execute("6+20=26", $false, $true, $false, $false, $false, $true, 1, $false, "86-62=24", 2) | Out-Null
$d.Content.Find.Execute("98-91=7", $false, $true, $false, $false, $false, $true, 1, $false, "80-56=24", 2) | Out-Null
$d.Content.Find.Execute("7+55=62", $false, $true, $false, $false, $false, $true, 1, $false, "83-5=78", 2) | Out-Null
$d.Content.Find.Execute("21+10=31", $false, $true, $false, $false, $false, $true, 1, $false, "60+17=77", 2) | Out-Null
$d.Content.Find.Execute("28+17=45", $false, $true, $false, $false, $false, $true, 1, $false, "9+25=34", 2) | Out-Null
$d.Content.Find.Execute("7+15=22", $false, $true, $false, $false, $false, $true, 1, $false, "84-10=74", 2) | Out-Null
$d.Content.Find.Execute("17+36=53", $false, $true, $false, $false, $false, $true, 1, $false, "64-46=18", 2) | Out-Null
$d.Content.Find.Execute("70+27=97", $false, $true, $false, $false, $false, $true, 1, $false, "82+16=98", 2) | Out-Null
$d.Content.Find.Execute("69-49=20", $false, $true, $false, $false, $false, $true, 1, $false, "0+43=43", 2) | Out-Null
$d.Content.Find.Execute("70-5=65", $false, $true, $false, $false, $false, $true, 1, $false, "25+5=30", 2) | Out-Null
$d.Content.Find.Execute("44-36=8", $false, $true, $false, $false, $false, $true, 1, $false, "57-10=47", 2) | Out-Null
$d.Content.Find.Execute("75-54=21", $false, $true, $false, $false, $false, $true, 1, $false, "0+29=29", 2) | Out-Null
$d.Content.Find.Execute("69-27=42", $false, $true, $false, $false, $false, $true, 1, $false, "93-92=1", 2) | Out-Null
$d.Content.Find.Execute("47-4=43", $false, $true, $false, $false, $false, $true, 1, $false, "40+44=84", 2) | Out-Null
$d.Content.Find.Execute("39+25=64", $false, $true, $false, $false, $false, $true, 1, $false, "60-25=35", 2) | Out-Null
$d.Content.Find.Execute("51-45=6", $false, $true, $false, $false, $false, $true, 1, $false, "78+3=81", 2) | Out-Null
$d.Content.Find.Execute("4+10=14", $false, $true, $false, $false, $false, $true, 1, $false, "38-0=38", 2) | Out-Null
$d.Content.Find.Execute("21+72=93", $false, $true, $false, $false, $false, $true, 1, $false, "40+51=91", 2) | Out-Null
$d.Content.Find.Execute("56+17=73", $false, $true, $false, $false, $false, $true, 1, $false, "71-66=5", 2) | Out-Null
$d.Content.Find.Execute("25+22=47", $false, $true, $false, $false, $false, $true, 1, $false, "75+3=78", 2) | Out-Null
$d.Content.Find.Execute("87-23=64", $false, $true, $false, $false, $false, $true, 1, $false, "58-0=58", 2) | Out-Null
$d.Content.Find.Execute("72-11=61", $false, $true, $false, $false, $false, $true, 1, $false, "80-70=10", 2) | Out-Null
$d.Content.Find.Execute("44+18=62", $false, $true, $false, $false, $false, $true, 1, $false, "29+50=79", 2) | Out-Null
$d.Content.Find.Execute("0+39=39", $false, $true, $false, $false, $false, $true, 1, $false, "63-57=6", 2) | Out-Null
$d.Content.Find.Execute("41-37=4", $false, $true, $false, $false, $false, $true, 1, $false, "46+51=97", 2) | Out-Null
$d.Content.Find.Execute("35-1=34", $false, $true, $false, $false, $false, $true, 1, $false, "6-3=3", 2) | Out-Null
$d.Content.Find.Execute("15+2=17", $false, $true, $false, $false, $false, $true, 1, $false, "84-4=80", 2) | Out-Null
$d.Content.Find.Execute("13+70=83", $false, $true, $false, $false, $false, $true, 1, $false, "59+28=87", 2) | Out-Null
$d.Content.Find.Execute("61+36=97", $false, $true, $false, $false, $false, $true, 1, $false, "21-15=6", 2) | Out-Null
$d.Content.Find.Execute("17+9=26", $false, $true, $false, $false, $false, $true, 1, $false, "66-35=31", 2) | Out-Null
$d.Content.Find.Execute("48-38=10", $false, $true, $false, $false, $false, $true, 1, $false, "55-34=21", 2) | Out-Null
$d.Content.Find.Execute("56-49=7", $false, $true, $false, $false, $false, $true, 1, $false, "81-54=27", 2) | Out-Null
$d.Content.Find.Execute("49+21=70", $false, $true, $false, $false, $false, $true, 1, $false, "60-30=30", 2) | Out-Null
$d.Content.Find.Execute("84-68=16", $false, $true, $false, $false, $false, $true, 1, $false, "7+67=74", 2) | Out-Null
$d.Content.Find.Execute("57-22=35", $false, $true, $false, $false, $false, $true, 1, $false, "51+30=81", 2) | Out-Null
$d.Content.Find.Execute("90-73=17", $false, $true, $false, $false, $false, $true, 1, $false, "63-4=59", 2) | Out-Null
$d.Content.Find.Execute("0+45=45", $false, $true, $false, $false, $false, $true, 1, $false, "82-24=58", 2) | Out-Null
$d.Content.Find.Execute("66-65=1", $false, $true, $false, $false, $false, $true, 1, $false, "67-16=51", 2) | Out-Null
$d.Content.Find.Execute("18+37=55", $false, $true, $false, $false, $false, $true, 1, $false, "47-22=25", 2) | Out-Null
$d.Content.Find.Execute("67-18=49", $false, $true, $false, $false, $false, $true, 1, $false, "6+5=11", 2) | Out-Null
$d.Content.Find.Execute("54+22=76", $false, $true, $false, $false, $false, $true, 1, $false, "1+0=1", 2) | Out-Null
$d.Content.Find.Execute("24+74=98", $false, $true, $false, $false, $false, $true, 1, $false, "15+45=60", 2) | Out-Null
$d.Content.Find.Execute("34-0=34", $false, $true, $false, $false, $false, $true, 1, $false, "93-9=84", 2) | Out-Null
$d.Content.Find.Execute("7+75=82", $false, $true, $false, $false, $false, $true, 1, $false, "80-22=58", 2) | Out-Null
$d.Content.Find.Execute("52+47=99", $false, $true, $false, $false, $false, $true, 1, $false, "81-72=9", 2) | Out-Null
$d.Content.Find.Execute("43-3=40", $false, $true, $false, $false, $false, $true, 1, $false, "48-9=39", 2) | Out-Null
$d.Content.Find.Execute("5+71=76", $false, $true, $false, $false, $false, $true, 1, $false, "75-70=5", 2) | Out-Null
$d.Content.Find.Execute("29+34=63", $false, $true, $false, $false, $false, $true, 1, $false, "32+5=37", 2) | Out-Null
$d.Content.Find.Execute("32-27=5", $false, $true, $false, $false, $false, $true, 1, $false, "67-40=27", 2) | Out-Null
$d.Content.Find.Execute("91-29=62", $false, $true, $false, $false, $false, $true, 1, $false, "71-32=39", 2) | Out-Null
